$d = $word.ActiveDocument

$newText = "Dates à utiliser pour la Campagne Bootes: 14-23 mai, 13-22 juin, 12-21 juillet"

# Two textual variants exist in the document: two paragraphs end with an
# extra (differently-formatted) trailing space run, two do not. Handle the
# "with trailing space" variant first so its search string (being longer)
# is not pre-empted by the shorter "without trailing space" variant.
$searchWithSpace = "Dates à utiliser pour la Campagne 2018 Persée:  Du 30 octobre au 8 novembre et du 29 novembre au 8 décembre "
$searchNoSpace   = "Dates à utiliser pour la Campagne 2018 Persée:  Du 30 octobre au 8 novembre et du 29 novembre au 8 décembre"

$maxIterations = 10

$i = 0
while ($i -lt $maxIterations) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchWithSpace, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $rng.Delete()
    $rng.InsertAfter($newText)
    $i = $i + 1
}

$i = 0
while ($i -lt $maxIterations) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchNoSpace, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $rng.Delete()
    $rng.InsertAfter($newText)
    $i = $i + 1
}
